$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.398180842399597
$ws.Range("B1").Value = 1.947030901908875
$ws.Range("C1").Value = 3.460329294204712
$ws.Range("D1").Value = 3.639624357223511
$ws.Range("E1").Value = 0.8988903164863586
